$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 2 (해체작업 row); row 3 (설치작업 row) shifts up to become row 2
$ws.Rows("2:2").Delete()

# Move the active selection, mirroring the author's final cursor position
$ws.Range("G6").Select()
